$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("training_schedule")

$ws.Range("B2").Value = "Welcome/introduction"

$ws.Range("B55").Value = "Data Monitor Training (in parallel)"
$ws.Range("B57").Value = "Data Monitor Training (in parallel)"
$ws.Range("B59").Value = "Data Monitor Training (in parallel)"

$ws.Range("B45").Value = "Supervisor Training (after)"
$ws.Range("B49").Value = "Supervisor Training (after)"
$ws.Range("B51").Value = "Supervisor Training (after)"

$ws.Activate()
$ws.Range("B46").Select()
